$d = $word.ActiveDocument

# 1) "donner" -> "données" (typo fix) and add "comme" before "inutilisables."
$d.Content.Find.Execute("de donner avant", $true, $false, $false, $false, $false, $true, 1, $false, "de données avant", 2)
$d.Content.Find.Execute("fonctions inutilisables.", $true, $false, $false, $false, $false, $true, 1, $false, "fonctions comme inutilisables.", 2)

# 2) "requêtes pour aussi grande échelle" -> "requêtes à grande échelle"
$d.Content.Find.Execute("requêtes pour aussi grande échelle", $true, $false, $false, $false, $false, $true, 1, $false, "requêtes à grande échelle", 2)

# 3) "graphe, avec prédécesseur" -> "graphe avec : prédécesseur"
$d.Content.Find.Execute("dictionnaire du graphe, avec prédécesseur", $true, $false, $false, $false, $false, $true, 1, $false, "dictionnaire du graphe avec : prédécesseur", 2)

# 4) "(utiliser une autre librairie)" -> "(en utilisant une autre librairie par exemple)"
$d.Content.Find.Execute("graphe (utiliser une autre librairie) ou alors", $true, $false, $false, $false, $false, $true, 1, $false, "graphe (en utilisant une autre librairie par exemple) ou alors", 2)

# 5) Rework the final sentence: split into two sentences and pluralize several words
$d.Content.Find.Execute("plus de connaissances, malgré le nombre d’information disponible sur l’internet, de grand problème mathématique reste très peu développé voir non résolu, comme", $true, $false, $false, $false, $false, $true, 1, $false, "plus de connaissances. Malgré le nombre d’informations disponibles sur Internet, de grands problèmes mathématiques restent très peu développés voir non résolus, comme", 2)
